$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Valor" column (B2:B10) held fractions that were displayed with a
# 0.0% number format (e.g. -0.989 shown as "-98.9%"). Rescale each value to
# a plain percentage-point number (x100) and swap the percentage format for
# a 1-decimal accounting-style numeric format.
$acctFormat = '_-* #,##0.0_-;\-* #,##0.0_-;_-* "-"??_-;_-@_-'

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Range("B$row")
    $cell.Value = $cell.Value2 * 100
    $cell.NumberFormat = $acctFormat
}

# Drop the leftover selection anchor (B3) from the author's editing session.
$ws.Range("A1").Select() | Out-Null
